# Update F-column (signup/participant count) values on sheets "展览" and "全部类型"
# to match data refreshed at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    2  = 2823
    4  = 367
    5  = 1574
    6  = 1159
    12 = 9519
    16 = 268
    19 = 680
    20 = 680
    21 = 1196
    22 = 1004
    23 = 2964
    24 = 2240
    25 = 1932
    29 = 1557
    30 = 301
    32 = 176
    33 = 220
    34 = 31
    35 = 342
    37 = 308
    39 = 25
    40 = 121
    41 = 1536
    42 = 129
    43 = 1479
    44 = 26
    47 = 363
    48 = 739
    50 = 313
}

$sheet4Updates = @{
    2  = 2823
    3  = 367
    4  = 1574
    6  = 1159
    8  = 9519
    14 = 268
    16 = 680
    17 = 680
    18 = 1196
    19 = 1004
    20 = 2964
    21 = 2240
    22 = 1932
    24 = 1557
    25 = 301
    27 = 176
    28 = 220
    29 = 31
    30 = 342
    32 = 308
    37 = 25
    38 = 121
    39 = 1536
    41 = 129
    42 = 1479
    43 = 26
    47 = 363
    48 = 739
    49 = 313
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
